$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- survey sheet: insert the new "dbs_preleve" question as row 24 -----
# (existing rows 24..34 shift down to 25..35)
$survey.Rows.Item(24).Insert()

$survey.Range("A24").Value2 = "select_one yes_no"
$survey.Range("B24").Value2 = "dbs_preleve"
$survey.Range("C24").Value2 = "Le DBS a-t-il été prélevé?"
$survey.Range("H24").Value2 = '${add_participant} = ''Oui'''
$survey.Range("J24").Value2 = "yes"

# match the look of the surrounding rows (vertical-center / wrap-text)
$survey.Range("A24").VerticalAlignment = -4108
$survey.Range("B24").VerticalAlignment = -4108
$survey.Range("B24").WrapText = $true
$survey.Range("C24").VerticalAlignment = -4108
$survey.Range("C24").WrapText = $true
$survey.Range("D24").VerticalAlignment = -4108
$survey.Range("D24").WrapText = $true
$survey.Range("H24").VerticalAlignment = -4108
$survey.Range("H24").WrapText = $true
$survey.Range("J24").VerticalAlignment = -4108
$survey.Range("J24").WrapText = $true

# --- settings sheet: bump form_id / form_title to the "V2" form -------
$settings.Range("B2").Value2 = "cg_oncho_oem_202303_2_biopsiev_2"
$settings.Range("A2").Value2 = "(2023 Mars) CEO - 2 Formulaire Biopsie V2"
$settings.Range("A3").Select()

# --- restore cursor / active sheet to match the authored file ---------
$survey.Activate()
$survey.Range("B25").Select()
